# edit.ps1 - applies the tracked-change edit described by the diff:
#  - Slide "問題の説明" (position 2, sldId 266) -> "作りたい関数"
#    plus several shapes nudged down, and a now-redundant caption shape removed.
#  - Slide "問題の解法" (positions 3 & 4, sldId 262 / 264) -> "関数の作り方"

$p = $ppt.ActivePresentation

function Get-ShapeById($slide, [int]$id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# ---- Slide 2 (sldId 266): "問題の説明" -> "作りたい関数" ----
$s2 = $p.Slides.Item(2)

$title2 = Get-ShapeById $s2 469
$title2.TextFrame.TextRange.Text = "作りたい関数"

# Nudge several shapes down (their Top position moves by +122944 EMU each;
# Left/Width/Height are unchanged). The point values below are chosen so
# that, after round-tripping through the COM Single-precision Top property,
# they land on the exact target EMU.
(Get-ShapeById $s2 43).Top  = 62.04189176377953   # 556375,664988 -> 556375,787932
(Get-ShapeById $s2 45).Top  = 72.41456992913386   # 2525126,796721 -> 2525126,919665
(Get-ShapeById $s2 46).Top  = 72.41456992913386   # 2896637,796721 -> 2896637,919665
(Get-ShapeById $s2 10).Top  = 185.95015748031497  # 4514975,2238623 -> 4514975,2361567
(Get-ShapeById $s2 11).Top  = 234.71259842519686  # 4082534,2857906 -> 4082534,2980850
(Get-ShapeById $s2 12).Top  = 72.41456992913386   # 4690661,796721 -> 4690661,919665
(Get-ShapeById $s2 15).Top  = 154.2240227480315   # 6910942,1835701 -> 6910942,1958645

# Shape 13 ("関数christmasTreeを定義する" caption) is removed entirely.
$del13 = Get-ShapeById $s2 13
if ($del13 -ne $null) {
    $del13.Delete()
}

# ---- Slide 3 (sldId 262): "問題の解法" -> "関数の作り方" ----
$s3 = $p.Slides.Item(3)
$title3 = Get-ShapeById $s3 469
$title3.TextFrame.TextRange.Text = "関数の作り方"

# ---- Slide 4 (sldId 264): "問題の解法" -> "関数の作り方" ----
$s4 = $p.Slides.Item(4)
$title4 = Get-ShapeById $s4 469
$title4.TextFrame.TextRange.Text = "関数の作り方"
